$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the displayed text for rows 2 and 3 (columns A and E), matching
# the shared-string reorder in the diff (hyperlink targets stay as-is).
$ws.Range("A2").Value = "World of Boats at Eyemouth Maritime Centre (Museum) ~ A World Class Boat and Dinghy Collection from all over the World. Visit Us at the Eyemouth Maritime Centre (Museum)"
$ws.Range("A3").Value = "the Silver Darlings"

$ws.Range("E2").Value = "https://web.archive.org/web/20101002230055/http://www.worldofboats.org/localarea/eyemouth-disaster"
$ws.Range("E3").Value = "http://sites.scran.ac.uk/secf_final/danger/links/link3.php"
